# marsh_alec.xlsx: add a "Save" column (H) — header + the single data row's value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", formatted like the other header cells (e.g. G1's
# bold/centered/bordered style). Copy G1's format into H1, then set the text so
# we don't clobber the formatting with a plain Value assignment.
$xlPasteFormats = -4122
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H1").Value = "Save"

# New data cell H2 = 0 (plain number, no special style — matches B2:G2).
$ws.Range("H2").Value = 0
